$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of column G (the "K" column) cell references to their new values.
# These values were regenerated (K instead of Strike#, with recomputed std/mean
# and s_vals) per the commit message.
$updates = @{
    "G2" = 0
    "G5" = 2
    "G6" = 1
    "G7" = 2
    "G8" = 1
    "G9" = 2
    "G10" = 0
    "G11" = 2
    "G12" = 3
    "G13" = 1
    "G14" = 2
    "G15" = 0
    "G16" = 0
    "G17" = 2
    "G18" = 1
    "G19" = 1
    "G20" = 0
    "G21" = 2
    "G22" = 1
    "G23" = 1
    "G24" = 1
    "G25" = 3
    "G26" = 0
    "G28" = 2
    "G29" = 0
    "G30" = 1
    "G31" = 2
    "G32" = 1
    "G33" = 2
    "G34" = 3
    "G35" = 2
    "G36" = 1
    "G37" = 1
    "G38" = 1
    "G39" = 1
    "G40" = 1
    "G41" = 0
    "G42" = 1
    "G43" = 1
    "G44" = 3
    "G45" = 0
    "G46" = 0
    "G47" = 1
    "G48" = 1
    "G49" = 3
    "G50" = 3
    "G51" = 2
    "G53" = 3
    "G54" = 3
    "G55" = 3
    "G56" = 2
    "G57" = 1
    "G58" = 0
    "G59" = 1
    "G60" = 1
    "G61" = 1
    "G62" = 1
    "G63" = 3
    "G64" = 1
    "G65" = 0
    "G66" = 2
    "G67" = 2
    "G68" = 1
    "G69" = 1
    "G70" = 1
    "G71" = 0
    "G72" = 1
    "G73" = 2
    "G75" = 0
    "G76" = 0
    "G77" = 2
    "G78" = 3
    "G79" = 2
    "G81" = 0
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
